# Update ticket-availability figures ("想去人数" / "最低票价") that changed
# between the previous scrape and the new one recorded at commit 456a3b4.
#
# Sheet "展览" (Exhibitions):
#   - G3: was the literal text "已售罄" (sold out) -> now a price, 70
#   - F5, F6, F7, F8, F9, F10, F11, F15, F21, F22, F24, F25: updated "want to go" counts
#
# Sheet "演出" (Shows):
#   - F2: updated "want to go" count
#
# Sheet "全部类型" (All types, mirrors 展览+演出 with an extra inserted row):
#   - G3: same "已售罄" -> 70 fix
#   - F5, F6, F7, F8, F9, F10, F11, F12, F16, F22, F23, F25, F26: updated counts

$wb = $excel.ActiveWorkbook

# ---- 展览 ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("G3").Value = 70
$ws1.Range("F5").Value = 512
$ws1.Range("F6").Value = 724
$ws1.Range("F7").Value = 1413
$ws1.Range("F8").Value = 217
$ws1.Range("F9").Value = 78
$ws1.Range("F10").Value = 120
$ws1.Range("F11").Value = 6051
$ws1.Range("F15").Value = 4851
$ws1.Range("F21").Value = 50
$ws1.Range("F22").Value = 4
$ws1.Range("F24").Value = 21
$ws1.Range("F25").Value = 3322

# ---- 演出 ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 60

# ---- 全部类型 ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G3").Value = 70
$ws4.Range("F5").Value = 60
$ws4.Range("F6").Value = 512
$ws4.Range("F7").Value = 724
$ws4.Range("F8").Value = 1413
$ws4.Range("F9").Value = 217
$ws4.Range("F10").Value = 78
$ws4.Range("F11").Value = 120
$ws4.Range("F12").Value = 6051
$ws4.Range("F16").Value = 4851
$ws4.Range("F22").Value = 50
$ws4.Range("F23").Value = 4
$ws4.Range("F25").Value = 21
$ws4.Range("F26").Value = 3322
